# Update "Design 1 Data" sheet: insert a new A/B row at 12 ("max_payload"/100000),
# shifting the existing general-parameter rows 12-40 down to 13-41, then refresh
# every recomputed value in the design/ferry/altitude data tables (columns D:K).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Design 1 Data")

# Shift columns A:B only (D:K stay put - they are a separate computed table)
for ($r = 40; $r -ge 12; $r--) {
    $srcA = $ws.Cells.Item($r, 1).Value2
    $srcB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 1).Value = $srcA
    $ws.Cells.Item($r + 1, 2).Value = $srcB
}

# New parameter inserted at row 12
$ws.Range("A12").Value = "max_payload"

# Updated parameter / recomputed values for "Design 1 Data"
$ws.Range("E2").Value = "51507919.18927524"
$ws.Range("H2").Value = "51507919.18927524"
$ws.Range("K2").Value = "51507919.18927524"
$ws.Range("E4").Value = "265905.7682338378"
$ws.Range("H4").Value = "74123.39690932246"
$ws.Range("K4").Value = "247105.4339477974"
$ws.Range("E5").Value = "2608535.586373949"
$ws.Range("H5").Value = "727150.5236804534"
$ws.Range("K5").Value = "2424104.307027892"
$ws.Range("E6").Value = "1289311.913317399"
$ws.Range("H6").Value = "-259886.5473271068"
$ws.Range("K6").Value = "1223994.993373831"
$ws.Range("B7").Value = "90000"
$ws.Range("E7").Value = "2172211.913317399"
$ws.Range("H7").Value = "623013.4526728932"
$ws.Range("K7").Value = "2106894.993373831"
$ws.Range("E8").Value = "1285142.663317399"
$ws.Range("H8").Value = "-264055.7973271068"
$ws.Range("K8").Value = "1219825.743373831"
$ws.Range("E9").Value = "436323.6730565503"
$ws.Range("H9").Value = "104137.0710075602"
$ws.Range("K9").Value = "317209.3136540609"
$ws.Range("E10").Value = "436323.6730565503"
$ws.Range("H10").Value = "104137.0710075602"
$ws.Range("K10").Value = "317209.3136540609"
$ws.Range("B12").Value = "100000"
$ws.Range("E12").Value = "782.6138047231788"
$ws.Range("H12").Value = "218.1569152626498"
$ws.Range("K12").Value = "727.2722528154185"
$ws.Range("B13").Value = "943600"
$ws.Range("B14").Value = "538000"
$ws.Range("E14").Value = "96.90905869256055"
$ws.Range("H14").Value = "51.16525171590381"
$ws.Range("K14").Value = "93.41984282680538"
$ws.Range("B15").Value = "90000"
$ws.Range("E15").Value = "8.075754891046712"
$ws.Range("H15").Value = "4.263770976325318"
$ws.Range("K15").Value = "7.784986902233783"
$ws.Range("B16").Value = "425"
$ws.Range("B17").Value = "115.749"
$ws.Range("E17").Value = "0.04127581109511973"
$ws.Range("H17").Value = "0.07817805768278223"
$ws.Range("K17").Value = "0.04281745589548629"
$ws.Range("B18").Value = "2.7e-05"
$ws.Range("E18").Value = "1.83527261507357"
$ws.Range("H18").Value = "1.49884578273208"
$ws.Range("K18").Value = "1.812794897626504"
$ws.Range("B19").Value = "2.1e-08"
$ws.Range("E19").Value = "0.1101440556657076"
$ws.Range("H19").Value = "0.1101440556657076"
$ws.Range("K19").Value = "0.1101440556657076"
$ws.Range("B20").Value = "0.85"
$ws.Range("B21").Value = "0.02"
$ws.Range("B22").Value = "0.85"
$ws.Range("E22").Value = "0.1162209426725255"
$ws.Range("K22").Value = "0.295726056294236"
$ws.Range("B23").Value = "12"
$ws.Range("H23").Value = "6601813.591169997"
$ws.Range("B24").Value = "5"
$ws.Range("E24").Value = "23682944.76363733"
$ws.Range("K24").Value = "22008489.63093534"
$ws.Range("B25").Value = "0.001"
$ws.Range("B26").Value = "0"
$ws.Range("B27").Value = "1.812794897626504"
$ws.Range("B28").Value = "8"
$ws.Range("B29").Value = "1.5"
$ws.Range("B30").Value = "1.6"
$ws.Range("B31").Value = "1.8"
$ws.Range("B32").Value = "77.16"
$ws.Range("B33").Value = "61.728"
$ws.Range("B34").Value = "51.44"
$ws.Range("B35").Value = "3048"
$ws.Range("B36").Value = "70"
$ws.Range("B37").Value = "4"
$ws.Range("B38").Value = "251.3274122871834"
$ws.Range("B39").Value = "9.81"

# Refresh recomputed OEW/ZFW/EW values on the other design sheets

$ws = $wb.Worksheets.Item("Design 2 Data")
$ws.Range("E6").Value = "1301166.277536666"
$ws.Range("H6").Value = "-237910.1813767925"
$ws.Range("K6").Value = "1461971.409645771"
$ws.Range("E8").Value = "1296997.027536666"
$ws.Range("H8").Value = "-242079.4313767925"
$ws.Range("K8").Value = "1457802.159645771"

$ws = $wb.Worksheets.Item("Design 3 Data")
$ws.Range("E6").Value = "1293060.394926933"
$ws.Range("H6").Value = "-242500.6570095802"
$ws.Range("K6").Value = "1235575.482640902"
$ws.Range("E8").Value = "1288891.144926933"
$ws.Range("H8").Value = "-246669.9070095802"
$ws.Range("K8").Value = "1231406.232640902"

$ws = $wb.Worksheets.Item("Design 4 Data")
$ws.Range("E6").Value = "1293396.414200417"
$ws.Range("H6").Value = "-242350.2442693433"
$ws.Range("K6").Value = "1235654.70814583"
$ws.Range("E8").Value = "1289227.164200417"
$ws.Range("H8").Value = "-246519.4942693433"
$ws.Range("K8").Value = "1231485.45814583"
